$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3545
$ws1.Range("F5").Value = 2203
$ws1.Range("F6").Value = 430
$ws1.Range("F12").Value = 1808
$ws1.Range("F13").Value = 135

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3545
$ws4.Range("F5").Value = 2203
$ws4.Range("F6").Value = 430
$ws4.Range("F15").Value = 1808
$ws4.Range("F16").Value = 135
